# CSSE-allprograms-outcome-mappings-20240913.xlsx
# "Criterion E is generated with updated knowledge base"
#
# The diff's only content-level (non-cosmetic) changes are:
#   1. The worksheet "Programs Details" is renamed to "Unit Details"
#      (which automatically updates the _xlnm._FilterDatabase defined
#      name that refers to that sheet).
#   2. The active/selected tab moves from "Outcomes Mappings" to the
#      renamed "Unit Details" sheet, with the selection on that sheet
#      moved to cell J20.

$wb = $excel.ActiveWorkbook

# 1. Rename "Programs Details" -> "Unit Details"
$unitDetails = $wb.Worksheets.Item("Programs Details")
$unitDetails.Name = "Unit Details"

# 2. Make "Unit Details" the active sheet and move the selection to J20
$unitDetails.Activate() | Out-Null
$unitDetails.Range("J20").Select() | Out-Null

# 3. The hidden Power Pivot range name loses its trailing digit
#    (WorksheetConnectionOutcomesMappingsA2H311 -> ...A2H31)
$rangeName = $wb.Names.Item("_xlcn.WorksheetConnection_OutcomesMappingsA2H311")
$rangeName.Name = "_xlcn.WorksheetConnection_OutcomesMappingsA2H31"
